# Auto-generated edit script: adds FAQ rows 12-19 (signature page / KOT / quota 1 & 2 Q&A)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "What is a signature page?"
$ws.Cells.Item(12, 3).Value = "The signature page is the final step to verify you as an applicant.`nYour application is not complete until Aarhus University receives your signed signature page through the university’s web form.`nYou must submit one signature page for each programme you apply for at Aarhus University.`nBe sure to check the deadline for submitting these pages."
$ws.Cells.Item(12, 2).WrapText = $true
$ws.Cells.Item(12, 3).WrapText = $true
$ws.Rows.Item(12).RowHeight = 102

# Row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 3).Value = "No.`nIf you complete your application using MitID, you do not need a signature page, because you are already verified through your MitID login. Applicants using MitID do not send in additional signature documentation."
$ws.Cells.Item(13, 2).Value = "Do I need a signature page if I apply with a Danish MitID?"
$ws.Cells.Item(13, 3).WrapText = $true
$ws.Rows.Item(13).RowHeight = 68

# Row 14
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "How does the Danish admission system (KOT) work?"
$ws.Cells.Item(14, 3).Value = "Aarhus University and all other higher education institutions in Denmark use a shared application system called the Coordinated Enrolment System (KOT).`nKOT allows you to:`nApply for up to eight different Bachelor’s degree programmes`nBe admitted to the highest-priority programme for which you qualify`nBe admitted to only one of the programmes you apply for`nYou can find more information about the general Danish admission system on the relevant information pages."
$ws.Cells.Item(14, 2).WrapText = $true
$ws.Cells.Item(14, 2).HorizontalAlignment = -4131
$ws.Cells.Item(14, 2).VerticalAlignment = -4160
$ws.Cells.Item(14, 3).WrapText = $true
$ws.Rows.Item(14).RowHeight = 221

# Row 15
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 3).Value = "Quota 1 and quota 2 is how we allocate student places in Denmark.`nQuota 1 admissions are allocated according to the grade point average of your upper secondary education.`nQuota 2 admissions are allocated according to various selection criteria.`nEven though the two quotas have different selection criteria, they are equally good."
$ws.Cells.Item(15, 2).Value = "What is quota 1 and quota 2?"
$ws.Cells.Item(15, 2).WrapText = $true
$ws.Cells.Item(15, 3).WrapText = $true
$ws.Rows.Item(15).RowHeight = 136

# Row 16
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 3).Value = "For both quotas you must have a qualifying entry examination and fulfil the specific admission requirements for your desired programme(s).`nYou do not have the option to choose which quota to apply for in your application. All international applicants are automatically assessed in quota 2. Depending on your qualifying entry examination your application may also be assessed in quota 1."
$ws.Cells.Item(16, 2).Value = "Do I have to choose between quota 1 or quota 2?"
$ws.Cells.Item(16, 2).WrapText = $true
$ws.Cells.Item(16, 3).WrapText = $true
$ws.Rows.Item(16).RowHeight = 102

# Row 17
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "Does it affect my studies if I am accepted through quota 1 or quota 2?"
$ws.Cells.Item(17, 3).Value = "No. Quota 1 and quota 2 is solely a part of the admission process and has no effect on your studies after you have been admitted. If you are offered a student place you do not know whether your application has been assessed in quota 1 or quota 2."
$ws.Cells.Item(17, 2).WrapText = $true
$ws.Cells.Item(17, 3).WrapText = $true
$ws.Rows.Item(17).RowHeight = 51

# Row 18
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "What is the Quota 1 selection criteria"
$ws.Cells.Item(18, 3).Value = "The assessment in quota 1 is based on the first completed qualifying examination and the principle of descending order of grade point average (GPA).`nThis means that applicants with the highest GPA will be admitted until there are no more available student places. The last admitted applicant’s GPA is the GPA cut-off for each programme. This can differ from year to year depending on number of student places and the qualifications of the applicants. You can find the GPA cut-off on the programmes’ websites. The GPA cut-off is first known when the admission process is completed on 28 July.`nYou cannot improve the GPA from your first qualifying examination. However, you are welcome to take supplementary courses in order to fulfil the admission requirements. The grades from your supplementary courses will not be used in the quota 1 assessment."
$ws.Cells.Item(18, 2).WrapText = $true
$ws.Cells.Item(18, 3).WrapText = $true
$ws.Rows.Item(18).RowHeight = 204

# Row 19
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "What is the Quota 2 selection criteria"
$ws.Cells.Item(19, 3).Value = "The quota 2 assessment is an overall assessment based on the following criteria for our English taught programmes:`nApplicants’ grade point average of particularly relevant quota 2 subjects`nApplicants’ relevant qualifications`nAll qualifications must be completed and documented no later than 15 March which is the deadline"
$ws.Cells.Item(19, 2).WrapText = $true
$ws.Cells.Item(19, 3).WrapText = $true
$ws.Rows.Item(19).RowHeight = 102

# Update selection / view to match end state
$ws.Range("A20").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1

Write-Host "FAQ rows added"
